$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: was "Test flight computer (clrnc)" -> now "Arduino Nano" ---
$ws.Range("B4").Value = "Arduino Nano"
$ws.Range("C4").Value = "Micro Controller "
$ws.Range("D4").Value = 6.59
$ws.Range("E4").Value = ""
$ws.Range("F4").Value = 8
$ws.Range("G4").Formula = "=D4*F4"

# --- Row 5: new item "Adafruit BNO055" ---
$ws.Range("B5").Value = "Adafruit BNO055"
$ws.Range("C5").Value = "Accelerometer"
$ws.Range("D5").Value = 34.95
$ws.Range("F5").Value = 5
$ws.Range("G5").Formula = "=D5*F5"

# --- Row 6: new item "SparkFun MS5803-14BA" ---
$ws.Range("B6").Value = "SparkFun MS5803-14BA"
$ws.Range("C6").Value = "Barometer"
$ws.Range("D6").Value = 59.95
$ws.Range("F6").Value = 5
$ws.Range("G6").Formula = "=D6*F6"

# --- Rows 7-9: blank rows that belong to the table (only G formula) ---
$ws.Range("G7").Formula = "=D7*F7"
$ws.Range("G8").Formula = "=D8*F8"
$ws.Range("G9").Formula = "=D9*F9"

# --- Row 10: Totals row (was row 5) ---
$ws.Range("B10").Value = "Total"
$ws.Range("D10").Formula = "=SUBTOTAL(101,Table1[Price per])"
$ws.Range("F10").Formula = "=SUBTOTAL(109,Table1[Quantity])"
$ws.Range("G10").Formula = "=SUBTOTAL(109,Table1[Price total])"

# Clear the old row 5 total-row formulas that are no longer part of the totals row
# (now overwritten above as row 10, row 5 itself holds the Adafruit BNO055 data)

# --- Resize the table to cover the new range, keeping the totals row ---
$ws.ListObjects.Item(1).Resize($ws.Range("B2:G10"))

# --- Selection matches the authored file ---
$ws.Range("B7").Select()
